# Update the cryptocurrency price/volume table with freshly scraped values.
# Row layout: A=rank index, B=Coin name, C=Link, D=Price, E=Volume(1h)
# D (Price) values are stored as plain text in the source sheet (they use
# "." as both decimal and thousands separators, e.g. "60.196.21"), so a
# leading apostrophe is used to force text entry and avoid Excel's
# automatic number coercion (which would strip things like a trailing
# "..00").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "'60.196.21"
$ws.Cells.Item(2, 5).Value = "  +0.16%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "'2.421.38"
$ws.Cells.Item(3, 5).Value = "  -0.12%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.03%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "'554.57"

# Row 6 - Solana
$ws.Cells.Item(6, 4).Value = "'137.28"
$ws.Cells.Item(6, 5).Value = "  -1.03%  "

# Row 7 - USDC
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 5).Value = "  +0.00%  "

# Row 8 - XRP
$ws.Cells.Item(8, 5).Value = "  +1.60%  "

# Row 9 - Dogecoin
$ws.Cells.Item(9, 5).Value = "  -1.50%  "

# Row 10 - Toncoin
$ws.Cells.Item(10, 4).Value = "'5.69"
$ws.Cells.Item(10, 5).Value = "  -1.23%  "

# Row 11 - TRON
$ws.Cells.Item(11, 5).Value = "  -0.11%  "

# Row 12 - Cardano
$ws.Cells.Item(12, 5).Value = "  -1.78%  "

# Row 13 - Avalanche
$ws.Cells.Item(13, 5).Value = "  -0.14%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = "'2.852.33"
$ws.Cells.Item(14, 5).Value = "  -0.16%  "

# Row 15 - WrappedBTC
$ws.Cells.Item(15, 4).Value = "'60.104.92"
$ws.Cells.Item(15, 5).Value = "  +0.12%  "

# Row 16 - ShibaInu
$ws.Cells.Item(16, 5).Value = "  -1.05%  "

# Row 17 - WrappedEther
$ws.Cells.Item(17, 4).Value = "'2.425.57"
$ws.Cells.Item(17, 5).Value = "  +0.14%  "

# Row 18 - Chainlink
$ws.Cells.Item(18, 4).Value = "'11.27"
$ws.Cells.Item(18, 5).Value = "  -0.97%  "

# Row 19 - Polkadot
$ws.Cells.Item(19, 4).Value = "'4.50"
$ws.Cells.Item(19, 5).Value = "  +2.34%  "

# Row 20 - BitcoinCash
$ws.Cells.Item(20, 4).Value = "'327.34"
$ws.Cells.Item(20, 5).Value = "  -1.54%  "

# Row 21 - Uniswap
$ws.Cells.Item(21, 5).Value = "  -0.65%  "

# Row 22 - Dai
$ws.Cells.Item(22, 5).Value = "  +0.14%  "

# Row 23 - Litecoin
$ws.Cells.Item(23, 4).Value = "'65.27"
$ws.Cells.Item(23, 5).Value = "  +0.25%  "

# Row 24 - Kaspa
$ws.Cells.Item(24, 5).Value = "  +4.08%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Cells.Item(25, 4).Value = "'8.73"
$ws.Cells.Item(25, 5).Value = "  +1.52%  "

# Row 26 - Binance-PegBSC-USD
$ws.Cells.Item(26, 4).Value = "'0.998"

# Row 27 - Fetch.AI
$ws.Cells.Item(27, 5).Value = "  +2.68%  "

# Row 28 - PEPE
$ws.Cells.Item(28, 4).Value = "'0.0₃0773"
$ws.Cells.Item(28, 5).Value = "  -2.08%  "

# Row 29 - PancakeSwap
$ws.Cells.Item(29, 4).Value = "'1.77"
$ws.Cells.Item(29, 5).Value = "  -0.80%  "

# Row 30 - Monero
$ws.Cells.Item(30, 4).Value = "'170.43"
$ws.Cells.Item(30, 5).Value = "  +0.88%  "

# Row 31 - Aptos
$ws.Cells.Item(31, 5).Value = "  -3.09%  "

# Row 32 - SuiNetwork
$ws.Cells.Item(32, 4).Value = "'1.08"
$ws.Cells.Item(32, 5).Value = "  +1.64%  "

# Row 33 - PolygonEcosystemToken
$ws.Cells.Item(33, 5).Value = "  -4.28%  "

# Row 34 - EthereumClassic
$ws.Cells.Item(34, 4).Value = "'18.54"
$ws.Cells.Item(34, 5).Value = "  -0.90%  "

# Row 36 - ImmutableX
$ws.Cells.Item(36, 5).Value = "  +1.60%  "

# Row 37 - FirstDigitalUSD
$ws.Cells.Item(37, 5).Value = "  +0.07%  "

# Row 38 - NEARProtocol
$ws.Cells.Item(38, 5).Value = "  -0.23%  "

# Row 39 - Bittensor
$ws.Cells.Item(39, 4).Value = "'328.38"
$ws.Cells.Item(39, 5).Value = "  +1.97%  "

# Row 40 - Stacks
$ws.Cells.Item(40, 5).Value = "  -0.88%  "

# Row 41 - Aave
$ws.Cells.Item(41, 4).Value = "'145.01"
$ws.Cells.Item(41, 5).Value = "  +3.38%  "

# Row 42 - Filecoin
$ws.Cells.Item(42, 5).Value = "  -1.08%  "

# Rows 43/44 swap: InjectiveProtocol now ranks above Stellar
$ws.Cells.Item(43, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(43, 4).Value = "'20.02"
$ws.Cells.Item(43, 5).Value = "  +2.68%  "

$ws.Cells.Item(44, 2).Value = "Stellar"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(44, 4).Value = "'0.0966"
$ws.Cells.Item(44, 5).Value = "  +0.48%  "

# Row 45 - Hedera
$ws.Cells.Item(45, 5).Value = "  -1.18%  "

# Row 47 - VeChain
$ws.Cells.Item(47, 5).Value = "  -1.52%  "

# Row 48 - WhiteBITCoin
$ws.Cells.Item(48, 5).Value = "  -0.08%  "

# Row 49 - dogwifhat
$ws.Cells.Item(49, 5).Value = "  -2.76%  "

# Row 50 - ZEEBU
$ws.Cells.Item(50, 5).Value = "  -0.56%  "

# Row 51 - BitgetToken
$ws.Cells.Item(51, 5).Value = "  -0.67%  "
